# Fruta / hortaliza, semanal
# Update the weekly market-price rows: dates, volumes, prices and origin/unit
# info get rotated between rows 2,3,4,7,8,9,11,14 to reflect the refreshed
# weekly data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44187
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 2800
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 2900
$ws.Range("S2").Value = 1450

# Row 3
$ws.Range("D3").Value = 44187
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 1400
$ws.Range("O3").Value = 1500
$ws.Range("P3").Value = 1446
$ws.Range("Q3").Value = '$/envase 1 kilo'
$ws.Range("R3").Value = "Provincia de Diguillín"
$ws.Range("S3").Value = 1446
$ws.Range("T3").Value = 1

# Row 4
$ws.Range("D4").Value = 44539
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3800
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3900
$ws.Range("S4").Value = 1950

# Row 7
$ws.Range("D7").Value = 44594
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 2500
$ws.Range("O7").Value = 2800
$ws.Range("P7").Value = 2650
$ws.Range("S7").Value = 1325

# Row 8
$ws.Range("D8").Value = 44174
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 3700
$ws.Range("O8").Value = 3800
$ws.Range("P8").Value = 3747
$ws.Range("Q8").Value = '$/bandeja 2 kilos'
$ws.Range("R8").Value = "Provincia de Linares"
$ws.Range("S8").Value = 1874
$ws.Range("T8").Value = 2

# Row 9
$ws.Range("D9").Value = 44931
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("R9").Value = "Provincia de Diguillín"
$ws.Range("S9").Value = 1500

# Row 11
$ws.Range("D11").Value = 44540
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 3500
$ws.Range("O11").Value = 3800
$ws.Range("P11").Value = 3650
$ws.Range("R11").Value = "Región del Maule"
$ws.Range("S11").Value = 1825

# Row 14
$ws.Range("D14").Value = 44932
$ws.Range("M14").Value = 60
